$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Golden")

$ws.Range("B4").Value = 809.687
$ws.Range("C4").Value = 11237183.211
$ws.Range("D4").Value = 23340243.403
$ws.Range("E4").Value = 6379505.299
$ws.Range("F4").Value = 4.449
$ws.Range("G4").Value = 4.697
$ws.Range("H4").Value = 5.365
$ws.Range("B5").Value = 154862177.941
$ws.Range("C5").Value = 213142558.437
$ws.Range("D5").Value = 269288404.216
$ws.Range("E5").Value = 39778983.974
$ws.Range("F5").Value = 37.047
$ws.Range("G5").Value = 37.459
$ws.Range("H5").Value = 37.788
$ws.Range("B6").Value = 372451257.346
$ws.Range("C6").Value = 423759402.813
$ws.Range("D6").Value = 514936531.617
$ws.Range("E6").Value = 41476770.343
$ws.Range("F6").Value = 101.056
$ws.Range("G6").Value = 102.133
$ws.Range("H6").Value = 106.064
$ws.Range("B10").Value = 0.007
$ws.Range("C10").Value = 7592.682
$ws.Range("D10").Value = 13009.57
$ws.Range("E10").Value = 4979.597
$ws.Range("F10").Value = 2.045
$ws.Range("G10").Value = 2.126
$ws.Range("H10").Value = 2.287
$ws.Range("B11").Value = 17658.604
$ws.Range("C11").Value = 56757.67
$ws.Range("D11").Value = 66604.288
$ws.Range("E11").Value = 13810.082
$ws.Range("F11").Value = 13.757
$ws.Range("G11").Value = 14.278
$ws.Range("H11").Value = 14.764
$ws.Range("B12").Value = 102100.83
$ws.Range("C12").Value = 116331.477
$ws.Range("D12").Value = 126348.038
$ws.Range("E12").Value = 7407.646
$ws.Range("F12").Value = 34.922
$ws.Range("G12").Value = 35.352
$ws.Range("H12").Value = 36.189
$ws.Range("B16").Value = 2.776
$ws.Range("C16").Value = 15.715
$ws.Range("D16").Value = 19.775
$ws.Range("E16").Value = 6.371
$ws.Range("F16").Value = 3.99
$ws.Range("G16").Value = 4.054
$ws.Range("H16").Value = 4.176
$ws.Range("B17").Value = 20.134
$ws.Range("C17").Value = 20.536
$ws.Range("D17").Value = 20.774
$ws.Range("E17").Value = 0.158
$ws.Range("F17").Value = 28.267
$ws.Range("G17").Value = 28.653
$ws.Range("H17").Value = 29.258
$ws.Range("B18").Value = 20.549
$ws.Range("C18").Value = 20.724
$ws.Range("D18").Value = 20.869
$ws.Range("E18").Value = 0.095
$ws.Range("F18").Value = 73.71299999999999
$ws.Range("G18").Value = 75.56100000000001
$ws.Range("H18").Value = 78.52500000000001
$ws.Range("B22").Value = 1.065
$ws.Range("C22").Value = 74.062
$ws.Range("D22").Value = 116.332
$ws.Range("E22").Value = 39.959
$ws.Range("F22").Value = 4.074
$ws.Range("G22").Value = 4.239
$ws.Range("H22").Value = 4.4
$ws.Range("B23").Value = 446.43
$ws.Range("C23").Value = 543.7809999999999
$ws.Range("D23").Value = 617.864
$ws.Range("E23").Value = 50.938
$ws.Range("F23").Value = 30.074
$ws.Range("G23").Value = 31.088
$ws.Range("H23").Value = 32.174
$ws.Range("B24").Value = 989.727
$ws.Range("C24").Value = 1056.827
$ws.Range("D24").Value = 1130.995
$ws.Range("E24").Value = 43.659
$ws.Range("F24").Value = 80.224
$ws.Range("G24").Value = 82.98699999999999
$ws.Range("H24").Value = 85.14100000000001
$ws.Range("B28").Value = 93.831
$ws.Range("C28").Value = 206250.984
$ws.Range("D28").Value = 364318.054
$ws.Range("E28").Value = 119098.041
$ws.Range("F28").Value = 3.78
$ws.Range("G28").Value = 3.872
$ws.Range("H28").Value = 3.971
$ws.Range("B29").Value = 1262383.423
$ws.Range("C29").Value = 1518041.319
$ws.Range("D29").Value = 1809166.501
$ws.Range("E29").Value = 179280.57
$ws.Range("F29").Value = 28.632
$ws.Range("G29").Value = 29.165
$ws.Range("H29").Value = 30.322
$ws.Range("B30").Value = 2741063.758
$ws.Range("C30").Value = 3172679.587
$ws.Range("D30").Value = 3389324.577
$ws.Range("E30").Value = 196937.608
$ws.Range("F30").Value = 76.726
$ws.Range("G30").Value = 78.09399999999999
$ws.Range("H30").Value = 80.235
$ws.Range("B34").Value = 12.3
$ws.Range("C34").Value = 3863847.826
$ws.Range("D34").Value = 18455086.861
$ws.Range("E34").Value = 5494331.53
$ws.Range("F34").Value = 7.965
$ws.Range("G34").Value = 8.321999999999999
$ws.Range("H34").Value = 8.760999999999999
$ws.Range("B35").Value = 208679432.699
$ws.Range("C35").Value = 441762958.29
$ws.Range("D35").Value = 615676602.46
$ws.Range("E35").Value = 105949769.312
$ws.Range("F35").Value = 65.405
$ws.Range("G35").Value = 67.39100000000001
$ws.Range("H35").Value = 70.60599999999999
$ws.Range("B36").Value = 817405757.696
$ws.Range("C36").Value = 1031807824.065
$ws.Range("D36").Value = 1221361116.583
$ws.Range("E36").Value = 121164612.685
$ws.Range("F36").Value = 179.077
$ws.Range("G36").Value = 183.083
$ws.Range("H36").Value = 189.299
